# Auto-generated PowerShell COM-interop script
$p = $ppt.ActivePresentation

# --- 1. Update the cached datetimeFigureOut field text (5/27/2021 -> 5/28/2021) ---
# on the slide master and every slide layout that carries a Date placeholder.
$newDate = "5/28/2021"
$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    foreach ($shp in $layout.Shapes) {
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Append 8 new slides (7-14), cloned from slide 6 so they inherit its
#        title gradient fill + body normAutofit/justify/150% line spacing formatting ---
$lastIndex = 6
for ($i = 1; $i -le 8; $i++) {
    $src = $p.Slides.Item($lastIndex)
    [void]$src.Duplicate()
    $lastIndex = $lastIndex + 1
}

# --- 3. Set title + body text for the new slides ---
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Scalability"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Scalability is the property of a system to handle a growing amount of work by adding resources to the system.`r`rIn computing, scalability is a characteristic of computers, networks, algorithms, protocols, and applications."

$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Types of Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Vertical Scaling (Scaling Up)`rHorizontal Scaling (Scaling Out)"

$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Vertical Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "`rVertical Scaling refers to scaling by adding more power (e.g. CPU, RAM) to an existing machine.`r"

$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Vertical Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Vertical scaling is limited to the capacity of one machine, scaling beyond that capacity can involve downtime and has an upper hard limit, i.e. the scale of the hardware on which you are currently running."

$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Horizontal Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "`rHorizontal Scaling means scaling by adding more machines to your pool of resources."

$s = $p.Slides.Item(12)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Horizontal Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "In theory, adding more machines to the existing pool means you are not limited to the capacity of a single unit, making it possible to scale with less downtime."

$s = $p.Slides.Item(13)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Types of Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "One of the fundamental differences between horizontal scaling and vertical scaling is that horizontal scaling requires breaking a sequential piece of logic into smaller pieces, so that they can be executed in parallel across multiple machines."

$s = $p.Slides.Item(14)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Types of Scaling"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "In many respects, vertical scaling is easier because the logic really doesn’t need to change. Rather, you’re just running the same code on higher-spec machines."

